# Se creo el primer escenario, con el tambien el step definition de dicho escenario
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# 1) Cell values for the new scenario row (row 4).
#    Order matters: it controls the order new entries are appended to
#    sharedStrings.xml (C4, D4, F4, E4 matches the target unique-string order).
# ---------------------------------------------------------------------------
$ws.Range("C4").Value = "En la pagina de myShopify me registro como un nuevo usuario "
$ws.Range("D4").Value = "Encontrarse en el formulario de registro."
$ws.Range("F4").Value = "El usuario este registrado y en su nueva cuenta "
$ws.Range("E4").Value = "1. ir a la opcion del formulario de registro. 2. llenar los datos entregados para el registro en el formulario. 3. darle submit. 4. verifico que si registro mi usuario mirando mi perfil."

# ---------------------------------------------------------------------------
# 2) Formatting for row 4: center/middle alignment everywhere, with wrap text
#    on the description/steps/result columns (C, E, F) so the long text is
#    readable; B (scenario name) and D (precondition) stay un-wrapped.
# ---------------------------------------------------------------------------
$ws.Range("C4").HorizontalAlignment = $xlCenter
$ws.Range("C4").VerticalAlignment = $xlCenter
$ws.Range("C4").WrapText = $true
$ws.Range("C4").Copy() | Out-Null
$ws.Range("E4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F4").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B4").HorizontalAlignment = $xlCenter
$ws.Range("B4").VerticalAlignment = $xlCenter
$ws.Range("B4").Copy() | Out-Null
$ws.Range("D4").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Row height / column widths for the now much denser row 4.
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 155.25

$ws.Columns.Item(2).ColumnWidth = 22.59245  # -> B width ~23.43
$ws.Columns.Item(3).ColumnWidth = 41.16667  # -> C width 42
$ws.Columns.Item(4).ColumnWidth = 37.8776   # -> D width ~38.71
$ws.Columns.Item(5).ColumnWidth = 31.59245  # -> E width ~32.43
$ws.Columns.Item(6).ColumnWidth = 26.73698  # -> F width ~27.57

# ---------------------------------------------------------------------------
# 4) Selection moves from B6 to H4.
# ---------------------------------------------------------------------------
$ws.Range("H4").Select() | Out-Null
